# Add the 'ubnt_configuration' automation script (TrainScheduling_ltrailways_searchCompany /
# CompanyManagement.searchCompany) as a new header+sample row pair at the bottom of the
# testData sheet, mirroring the existing row-pair layout used for every other script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41: header row (param-name row), formatted like the other header rows (e.g. row 38) ---
$ws.Range("A38:D38").Copy()
$ws.Range("A41:D41").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(41).RowHeight = 15.95

$ws.Cells.Item(41, 1).Value = "TrainScheduling_ltrailways_searchCompany"
$ws.Cells.Item(41, 2).Value = "1"
$ws.Cells.Item(41, 3).Value = "CompanyManagement.searchCompany"
$ws.Cells.Item(41, 4).Value = "companyName"

# --- Row 42: sample-value row, formatted like the other sample rows (e.g. row 39) ---
$ws.Range("A39:D39").Copy()
$ws.Range("A42:D42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F39").Copy()
$ws.Range("F42").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows.Item(42).RowHeight = 15.95

$ws.Cells.Item(42, 1).Value = "TrainScheduling_ltrailways_searchCompany"
$ws.Cells.Item(42, 2).Value = "1"
$ws.Cells.Item(42, 3).Value = "CompanyManagement.searchCompany"
$ws.Cells.Item(42, 4).Value = "4"

# Match the author's final view/selection state (scrolled down, D42 selected).
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("D42").Select()
